# The deck originally had 3 slides:
#   1) the Salsa20 "Table 4 / x_i" overview diagram slide
#   2) the k0 / Table 6 keystream slide
#   3) the Salsa20(x) block-diagram slide
#
# The authored change removes slide 1 entirely, leaving the other two
# slides (which shift up to positions 1 and 2, keeping their own content
# untouched).
$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
